$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.32"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.84"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05997"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.390"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8185"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9552"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1422"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07433"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03315"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03049"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09410"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001596"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04825"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005911"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005543"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004148"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009864"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.672"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.418"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002901"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03996"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006412"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1073"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002901"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005805"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005112"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8602"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.004763"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
